# Apply updated "Weight" (D) and "Percent Change" (E) values for rows 2-30,
# and bump the "as of" date in the confidential disclaimer text (A33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected (locked cells), so unprotect for the duration of
# the edit and restore protection afterwards.
$ws.Unprotect()

# New values per row, taken from the diff: row -> @(D, E)
# D30 is unchanged (stays 1), so it is omitted from the D-updates.
$values = @{
    2  = @(0.01937340688741203,  -0.01247366203118416)
    3  = @(0.0179613146396248,    0.006313945224045936)
    4  = @(0.07560935350429848,   0.007032271490440944)
    5  = @(0.05458605643143968,   0.001174755803076133)
    6  = @(0.07134116255582018,   0.02410052183466083)
    7  = @(0.01972283032389801,  -0.0003537318712416582)
    8  = @(0.0337254475855939,   -0.02508780732563975)
    9  = @(0.02849300232810415,   0.01444111027756945)
    10 = @(0.02368449646509358,   0.02253725917848048)
    11 = @(0.02574808381174953,  -0.0008647526807331607)
    12 = @(0.02613179646396545,  -0.01266721576869534)
    13 = @(0.04235608771626699,   0.006844346317615546)
    14 = @(0.02363669937650034,  -0.0061543871988744)
    15 = @(0.04020180465181454,   0.006793879600637931)
    16 = @(0.03003572015962446,  -0.003405075489881226)
    17 = @(0.04556576612201736,   0.01847749602564441)
    18 = @(0.1167828343732911,    0.00544522741832143)
    19 = @(0.02878334737558359,  -0.0003094250881860994)
    20 = @(0.02404698245995287,  -0.02017283950617288)
    21 = @(0.0245565706746752,    0.01389686459170791)
    22 = @(0.01341242712428371,  -0.001460871874896186)
    23 = @(0.01471393294038783,   0.007566204287515976)
    24 = @(0.03049736284759023,   0.007373877491421377)
    25 = @(0.0108665641849587,   -0.006515859355790687)
    26 = @(0.03727638532884694,  -0.009079180006689946)
    27 = @(0.02331251042778098,   0.004100552683187875)
    28 = @(0.05396098332564424,   0.01501959975242428)
    29 = @(0.04361706991378118,  -0.001252382248843054)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("D$row").Value = $pair[0]
    $ws.Range("E$row").Value = $pair[1]
}

# Row 30 (Total): D30 is unchanged, only E30 changes.
$ws.Range("E30").Value = 0.004325042765034981

# Update the disclaimer text date from 2021-07-13 to 2021-07-14.
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

$ws.Range("A33").Value = $newText

# Restore sheet protection to its prior (locked/protected) state.
$ws.Protect()
